# Automatic update of files.
#
# The source export re-ran and (a) bumped the "Taxonsorteringsordning" sort
# key for every "Garnlav" (Alectoria sarmentosa) observation from 79243 to
# 79244, and (b) re-ordered a handful of observation rows (12/13/14 and
# 18/19) that share the same location/date group, without touching rows
# that weren't part of that re-order (10, 11, 15, 16, 17, 20).
#
# Below we reproduce the resulting cell values directly, row by row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: sort key bump only ---
$ws.Range("B10").Value = 79244

# --- Row 12 (now holds what used to be row 14's observation) ---
$ws.Range("A12").Value = 131187835
$ws.Range("B12").Value = 57073
$ws.Range("D12").Value = "LC"
$ws.Range("E12").Value = 100138
$ws.Range("F12").Value = "Tjäder"
$ws.Range("G12").Value = "Tetrao urogallus"
$ws.Range("H12").Value = "Linnaeus, 1758"
$ws.Range("J12").ClearContents()
$ws.Range("L12").ClearContents()
$ws.Range("M12").Value = "äldre spår"
$ws.Range("Q12").Value = 511382
$ws.Range("R12").Value = 6697458
$ws.Range("S12").Value = 25
$ws.Range("Z12").ClearContents()
$ws.Range("AB12").ClearContents()
$ws.Range("AC12").Value = "Betad tallkrona."
$ws.Range("AF12").ClearContents()
$ws.Range("AW12").Value = "Anna-Lena Thommson"
$ws.Range("AX12").Value = "Anna-Lena Thommson"

# --- Row 13 (now holds what used to be row 12's observation) ---
$ws.Range("A13").Value = 131191949
$ws.Range("B13").Value = 79244
$ws.Range("P13").Value = "Svartå, Dlr"
$ws.Range("Q13").Value = 511393
$ws.Range("R13").Value = 6697824
$ws.Range("S13").Value = 10
$ws.Range("Z13").Value = "10:33"
$ws.Range("AB13").Value = "10:33"
$ws.Range("AC13").ClearContents()
$ws.Range("AW13").Value = "Lars-Erik Nilsson"
$ws.Range("AX13").Value = "Lars-Erik Nilsson, Anna-Lena Thommson"

# --- Row 14 (now holds what used to be row 13's observation) ---
$ws.Range("A14").Value = 131187780
$ws.Range("B14").Value = 79244
$ws.Range("D14").Value = "NT"
$ws.Range("E14").Value = 6425
$ws.Range("F14").Value = "Garnlav"
$ws.Range("G14").Value = "Alectoria sarmentosa"
$ws.Range("H14").Value = "(Ach.) Ach."
$ws.Range("J14").ClearContents()
$ws.Range("L14").ClearContents()
$ws.Range("M14").ClearContents()
$ws.Range("P14").Value = "Svatå, Dlr"
$ws.Range("Q14").Value = 511335
$ws.Range("R14").Value = 6697864
$ws.Range("AC14").Value = "På gran."
$ws.Range("AF14").ClearContents()

# --- Row 17: sort key bump only ---
$ws.Range("B17").Value = 79244

# --- Row 18 (now holds what used to be row 19's observation) ---
$ws.Range("A18").Value = 131187762
$ws.Range("B18").Value = 79244
$ws.Range("D18").Value = "NT"
$ws.Range("E18").Value = 6425
$ws.Range("F18").Value = "Garnlav"
$ws.Range("G18").Value = "Alectoria sarmentosa"
$ws.Range("H18").Value = "(Ach.) Ach."
$ws.Range("J18").ClearContents()
$ws.Range("L18").ClearContents()
$ws.Range("M18").ClearContents()
$ws.Range("P18").Value = "Svartå, Dlr"
$ws.Range("Q18").Value = 511511
$ws.Range("R18").Value = 6697866
$ws.Range("AC18").Value = "På äldre tall."
$ws.Range("AF18").ClearContents()

# --- Row 19 (now holds what used to be row 18's observation) ---
$ws.Range("A19").Value = 131187791
$ws.Range("B19").Value = 57073
$ws.Range("D19").Value = "LC"
$ws.Range("E19").Value = 100138
$ws.Range("F19").Value = "Tjäder"
$ws.Range("G19").Value = "Tetrao urogallus"
$ws.Range("H19").Value = "Linnaeus, 1758"
$ws.Range("J19").ClearContents()
$ws.Range("L19").ClearContents()
$ws.Range("M19").Value = "färsk spillning"
$ws.Range("P19").Value = "Svatå, Dlr"
$ws.Range("Q19").Value = 511301
$ws.Range("R19").Value = 6697864
$ws.Range("AC19").ClearContents()
$ws.Range("AF19").ClearContents()

# --- Row 20: sort key bump only ---
$ws.Range("B20").Value = 79244

Write-Output "Done applying row updates."
